$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C header "Value" - reuse B1's formatting (bold + border), matching style index s="1"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Value"

# Mirror existing column B values into the new column C for rows 2-5
$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = $ws.Range("B3").Value2
$ws.Range("C4").Value = $ws.Range("B4").Value2
$ws.Range("C5").Value = $ws.Range("B5").Value2

# New rows 6-8 (column B left blank on these rows, column C gets "NA")
$ws.Range("A6").Value = "old platform"
$ws.Range("B6").Borders.LineStyle = 0
$ws.Range("C6").Value = "NA"

$ws.Range("A7").Value = "new platform"
$ws.Range("B7").Borders.LineStyle = 0
$ws.Range("C7").Value = "NA"

$ws.Range("A8").Value = "service"
$ws.Range("B8").Borders.LineStyle = 0
$ws.Range("C8").Value = "NA"
